# Add 2022-Q3 data:
#  - insert a new worksheet "2022-Q3" right after "总计", shifting the
#    existing quarterly sheets (2021-Q2, 2021-Q1, 2020-Q4) one tab to the right
#  - populate the new sheet with the fund-holdings table for 2022-Q3
#  - insert a matching summary row into the "总计" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (position 2)
# ---------------------------------------------------------------------------
# NOTE: grab the "Before" anchor sheet first, insert, *then* re-resolve any
# other sheet handles by name -- Worksheets.Add() shifts tab positions, and
# a handle obtained beforehand can silently end up pointing at the new sheet
# (or the wrong neighbour) if it was positional under the hood.
$totalSheet = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Re-resolve by name *after* the insert/rename above.
$q2Sheet = $wb.Worksheets.Item("2021-Q2")

# Reuse the existing "2021-Q2" sheet's header / index-column formatting so
# the new sheet matches the look of its siblings (bold, centered, bordered).
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)

# ---- headers ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- data rows (text-like numeric fields keep a leading apostrophe so they
#      are stored as text, matching the source data's formatting) ----
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'159617"
$newSheet.Range("C2").Value = "华夏中证智选500价值稳健策略ETF"
$newSheet.Range("D2").Value = "'2.93"
$newSheet.Range("E2").Value = "'97.05"
$newSheet.Range("F2").Value = "'1.43"
$newSheet.Range("G2").Value = "'0.0419"
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'006729"
$newSheet.Range("C3").Value = "万家中证500指数增强A"
$newSheet.Range("D3").Value = "'2.43"
$newSheet.Range("E3").Value = "'94.01"
$newSheet.Range("F3").Value = "'1.22"
$newSheet.Range("G3").Value = "'0.0296"
$newSheet.Range("H3").Value = 7

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'006730"
$newSheet.Range("C4").Value = "万家中证500指数增强C"
$newSheet.Range("D4").Value = "'1.66"
$newSheet.Range("E4").Value = "'94.01"
$newSheet.Range("F4").Value = "'1.22"
$newSheet.Range("G4").Value = "'0.0203"
$newSheet.Range("H4").Value = 7

# ---------------------------------------------------------------------------
# 2. Insert a new summary row for 2022-Q3 into the "总计" sheet, shifting the
#    existing rows down by one.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")

# Preserve the styled index column ("A") formatting for the new last row
# (row 5) by copying it from the row above before shuffling values down.
$ws1.Range("A4").Copy()
$ws1.Range("A5").PasteSpecial(-4122)

# Shift existing data rows 2..4 down to 3..5 (bottom-up so we never
# overwrite a row before reading it).
for ($r = 4; $r -ge 2; $r--) {
    $newR = $r + 1
    $ws1.Cells.Item($newR, 2).Value = $ws1.Cells.Item($r, 2).Value2
    $ws1.Cells.Item($newR, 3).Value = $ws1.Cells.Item($r, 3).Value2
    $ws1.Cells.Item($newR, 4).Value = $ws1.Cells.Item($r, 4).Value2
}

# New row 2: 2022-Q3 summary data
$ws1.Cells.Item(2, 1).Value = 0
$ws1.Cells.Item(2, 2).Value = "2022-Q3"
$ws1.Cells.Item(2, 3).Value = 3
$ws1.Cells.Item(2, 4).Value = 0.09

# Renumber the index column for the rows that shifted down
$ws1.Cells.Item(3, 1).Value = 1
$ws1.Cells.Item(4, 1).Value = 2
$ws1.Cells.Item(5, 1).Value = 3
